$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6889.294
$ws.Range("I76").Value = 6112.3
$ws.Range("K76").Value = 6112.3
$ws.Range("M76").Value = -5797.3

$ws.Range("H79").Value = 6889.294
$ws.Range("I79").Value = 6112.3
$ws.Range("K79").Value = 6112.3
$ws.Range("M79").Value = -5020.3

$ws.Range("H86").Value = 4678.4443
$ws.Range("I86").Value = 3837.7273
$ws.Range("J86").Value = 5999.5713
$ws.Range("K86").Value = 3837.7273
$ws.Range("L86").Value = 5999.5713
$ws.Range("M86").Value = -2714.7273
$ws.Range("N86").Value = -8245.5713

$ws.Range("H89").Value = 4678.4443
$ws.Range("I89").Value = 3837.7273
$ws.Range("J89").Value = 5999.5713
$ws.Range("K89").Value = 19188.6365
$ws.Range("L89").Value = 29997.8565
$ws.Range("M89").Value = -13572.6365
$ws.Range("N89").Value = -41229.85649999999

$ws.Range("H107").Value = 895.2353000000001
$ws.Range("J107").Value = 610.5
$ws.Range("L107").Value = 610.5
$ws.Range("N107").Value = -4450.5

$ws.Range("H132").Value = 3050.6206
$ws.Range("I132").Value = 3520.348
$ws.Range("K132").Value = 10561.044
$ws.Range("M132").Value = -8031.044

$ws.Range("H137").Value = 4221.5107
$ws.Range("I137").Value = 1444.7297
$ws.Range("J137").Value = 14495.6
$ws.Range("K137").Value = 4334.189100000001
$ws.Range("L137").Value = 43486.8
$ws.Range("M137").Value = -1784.189100000001
$ws.Range("N137").Value = -48586.8

$ws.Range("H141").Value = 5955.273
$ws.Range("I141").Value = 6489.222
$ws.Range("J141").Value = 3552.5
$ws.Range("K141").Value = 19467.666
$ws.Range("L141").Value = 10657.5
$ws.Range("M141").Value = -14287.666
$ws.Range("N141").Value = -21017.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 41025.668
$ws.Range("I60").Value = 41025.668
$ws.Range("K60").Value = 41025.668
$ws.Range("M60").Value = -40292.668

$ws.Range("H61").Value = 37038996
$ws.Range("I61").Value = 50001736
$ws.Range("J61").Value = 2599.7144
$ws.Range("K61").Value = 50001736
$ws.Range("L61").Value = 2599.7144
$ws.Range("M61").Value = -50001524
$ws.Range("N61").Value = -3023.7144

$ws.Range("H76").Value = 103146.25
$ws.Range("J76").Value = 103146.25
$ws.Range("L76").Value = 103146.25
$ws.Range("N76").Value = -103822.25

$ws.Range("H79").Value = 103146.25
$ws.Range("J79").Value = 103146.25
$ws.Range("L79").Value = 103146.25
$ws.Range("N79").Value = -105486.25

$ws.Range("H110").Value = 2392.0293
$ws.Range("I110").Value = 2555.4167
$ws.Range("K110").Value = 2555.4167
$ws.Range("M110").Value = -510.4167000000002

$ws.Range("H136").Value = 37038996
$ws.Range("I136").Value = 50001736
$ws.Range("J136").Value = 2599.7144
$ws.Range("K136").Value = 150005208
$ws.Range("L136").Value = 7799.1432
$ws.Range("M136").Value = -150002658
$ws.Range("N136").Value = -12899.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 760.7727
$ws.Range("I16").Value = 770.4706
$ws.Range("J16").Value = 727.8
$ws.Range("K16").Value = 770.4706
$ws.Range("L16").Value = 727.8
$ws.Range("M16").Value = -483.4706
$ws.Range("N16").Value = -1301.8

$ws.Range("H21").Value = 1500
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H31").Value = 48081370
$ws.Range("I31").Value = 3087.6667
$ws.Range("J31").Value = 156257500
$ws.Range("K31").Value = 3087.6667
$ws.Range("L31").Value = 156257500
$ws.Range("M31").Value = -2792.6667
$ws.Range("N31").Value = -156258090

$ws.Range("H34").Value = 48081370
$ws.Range("I34").Value = 3087.6667
$ws.Range("J34").Value = 156257500
$ws.Range("K34").Value = 3087.6667
$ws.Range("L34").Value = 156257500
$ws.Range("M34").Value = -2885.6667
$ws.Range("N34").Value = -156257904

$ws.Range("H99").Value = 9480.531000000001
$ws.Range("I99").Value = 4319.8887
$ws.Range("J99").Value = 11499.913
$ws.Range("K99").Value = 4319.8887
$ws.Range("L99").Value = 11499.913
$ws.Range("M99").Value = -2821.8887
$ws.Range("N99").Value = -14495.913

$ws.Range("H107").Value = 790.7143
$ws.Range("I107").Value = 679.9091
$ws.Range("K107").Value = 679.9091
$ws.Range("M107").Value = 1240.0909

$ws.Range("H113").Value = 760.7727
$ws.Range("I113").Value = 770.4706
$ws.Range("J113").Value = 727.8
$ws.Range("K113").Value = 770.4706
$ws.Range("L113").Value = 727.8
$ws.Range("M113").Value = 1399.5294
$ws.Range("N113").Value = -5067.8

$ws.Range("H126").Value = 9480.531000000001
$ws.Range("I126").Value = 4319.8887
$ws.Range("J126").Value = 11499.913
$ws.Range("K126").Value = 12959.6661
$ws.Range("L126").Value = 34499.739
$ws.Range("M126").Value = -10489.6661
$ws.Range("N126").Value = -39439.739

$ws.Range("H132").Value = 2236
$ws.Range("I132").Value = 2130.5908
$ws.Range("K132").Value = 6391.7724
$ws.Range("M132").Value = -3861.7724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5584
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 6600.8
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 19802.4
$ws.Range("N5").Value = -20026.4
$ws.Range("M5").Value = -1388

$ws.Range("H92").Value = 1037.6
$ws.Range("J92").Value = 1994
$ws.Range("L92").Value = 5982
$ws.Range("N92").Value = -8478

$ws.Range("H131").Value = 5953617.5
$ws.Range("J131").Value = 1857.1818
$ws.Range("L131").Value = 5571.5454
$ws.Range("N131").Value = -15651.5454

$ws.Range("H135").Value = 5584
$ws.Range("I135").Value = 500
$ws.Range("J135").Value = 6600.8
$ws.Range("K135").Value = 4500
$ws.Range("L135").Value = 59407.2
$ws.Range("N135").Value = -64477.2
$ws.Range("M135").Value = -1965

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1351634
$ws.Range("I2").Value = 3571483
$ws.Range("J2").Value = 421.52173
$ws.Range("K2").Value = 3571483
$ws.Range("L2").Value = 421.52173
$ws.Range("M2").Value = -3571370
$ws.Range("N2").Value = -647.5217299999999

$ws.Range("H131").Value = 94981.664
$ws.Range("J131").Value = 94981.664
$ws.Range("L131").Value = 94981.664
$ws.Range("N131").Value = -105061.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1976.4445
$ws.Range("I46").Value = 1025.7368
$ws.Range("J46").Value = 4234.375
$ws.Range("K46").Value = 1025.7368
$ws.Range("L46").Value = 4234.375
$ws.Range("M46").Value = -837.7367999999999
$ws.Range("N46").Value = -4610.375

$ws.Range("H132").Value = 117656904
$ws.Range("I132").Value = 3338
$ws.Range("J132").Value = 222237860
$ws.Range("K132").Value = 10014.0
$ws.Range("L132").Value = 666713580
$ws.Range("M132").Value = -7484
$ws.Range("N132").Value = -666718640
